# Daily attendance processing - 2025-10-23 17:19:43
# Swap the order of the two comma-separated entries in the "Recorded By"
# column (G) for the specific rows where the value changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,4,6,10,11,12,13,14,15,17,18,19,20,21,22,30,31,33,37,38,39,40,41,42,44,45,46,47,48,49,57,58,60,64,65,66,67,68,69,71,72,73,74,75,76,86,87,88,89,90,93,95,96,97,99,102,112,113,114,115,116,119,121,122,123,125,128,138,139,140,141,142,145,147,148,149,151,154)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    $parts = $val -split ', ', 2
    if ($parts.Count -eq 2) {
        $cell.Value2 = $parts[1] + ", " + $parts[0]
    }
}
